$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values for the affected rows, as part of the
# "repull data, push all data, mean calculation" update.
$ws.Range("F2").Value = -4
$ws.Range("F4").Value = -3
$ws.Range("F5").Value = -3
$ws.Range("F7").Value = -2
